# Auto-generated edit script applying scheduled market-data refresh
# to the Balmung_Profits workbook (FFXIV leve profit calculator).
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC!row15 (hunk 0: @@ -1373,22 +1373,22 @@)
$ws_ALC.Range("H15").Value = 406.20587
$ws_ALC.Range("I15").Value = 406.20587
$ws_ALC.Range("K15").Value = 1218.61761
$ws_ALC.Range("M15").Value = -1049.61761

# ALC!row132 (hunk 1: @@ -7211,22 +7211,22 @@)
$ws_ALC.Range("H132").Value = 2604.682
$ws_ALC.Range("I132").Value = 2565.45
$ws_ALC.Range("K132").Value = 7696.349999999999
$ws_ALC.Range("M132").Value = -5166.349999999999

# ALC!row137 (hunk 2: @@ -7462,22 +7462,22 @@)
$ws_ALC.Range("H137").Value = 33334000
$ws_ALC.Range("J137").Value = 33334000
$ws_ALC.Range("L137").Value = 100002000
$ws_ALC.Range("N137").Value = -100007100

# ALC!row138 (hunk 3: @@ -7511,19 +7511,22 @@)
$ws_ALC.Range("H138").Value = 5298.2856
$ws_ALC.Range("I138").Value = 19532.834
$ws_ALC.Range("K138").Value = 58598.50199999999
$ws_ALC.Range("M138").Value = -53458.50199999999

# ALC!row141 (hunk 4: @@ -7655,22 +7655,22 @@)
$ws_ALC.Range("H141").Value = 2079.0667
$ws_ALC.Range("I141").Value = 2120.4285
$ws_ALC.Range("K141").Value = 6361.2855
$ws_ALC.Range("M141").Value = -1181.2855

# ARM!row4 (hunk 5: @@ -7905,22 +7905,22 @@)
$ws_ARM.Range("H4").Value = 767.7
$ws_ARM.Range("I4").Value = 754.5714
$ws_ARM.Range("K4").Value = 754.5714
$ws_ARM.Range("M4").Value = -638.5714

# ARM!row14 (hunk 6: @@ -8401,25 +8401,25 @@)
$ws_ARM.Range("H14").Value = 12542.5
$ws_ARM.Range("J14").Value = 16582.334
$ws_ARM.Range("L14").Value = 16582.334
$ws_ARM.Range("N14").Value = -16932.334

# ARM!row17 (hunk 7: @@ -8551,19 +8551,22 @@)
$ws_ARM.Range("H17").Value = 5655
$ws_ARM.Range("J17").Value = 5655
$ws_ARM.Range("L17").Value = 5655
$ws_ARM.Range("N17").Value = -6001

# ARM!row25 (hunk 8: @@ -8937,25 +8940,25 @@)
$ws_ARM.Range("H25").Value = 17646.111
$ws_ARM.Range("J25").Value = 21002.285
$ws_ARM.Range("L25").Value = 21002.285
$ws_ARM.Range("N25").Value = -21806.285

# ARM!row32 (hunk 9: @@ -9283,22 +9286,22 @@)
$ws_ARM.Range("H32").Value = 171087.38
$ws_ARM.Range("I32").Value = 189755.22
$ws_ARM.Range("K32").Value = 189755.22
$ws_ARM.Range("M32").Value = -189468.22

# ARM!row34 (hunk 10: @@ -9387,23 +9390,26 @@)
$ws_ARM.Range("H34").Value = 166332.67
$ws_ARM.Range("J34").Value = 53332.332
$ws_ARM.Range("L34").Value = 53332.332
$ws_ARM.Range("N34").Value = -53874.332

# ARM!row122 (hunk 11: @@ -13657,22 +13663,22 @@)
$ws_ARM.Range("H122").Value = 1084.9231
$ws_ARM.Range("I122").Value = 1084.9231
$ws_ARM.Range("K122").Value = 3254.7693
$ws_ARM.Range("M122").Value = -804.7692999999999

# BSM!row16 (hunk 12: @@ -15408,19 +15414,22 @@)
$ws_BSM.Range("H16").Value = 2500
$ws_BSM.Range("J16").Value = 2500
$ws_BSM.Range("L16").Value = 2500
$ws_BSM.Range("N16").Value = -2840

# BSM!row86 (hunk 13: @@ -18784,22 +18793,22 @@)
$ws_BSM.Range("H86").Value = 3449.25
$ws_BSM.Range("I86").Value = 1615.8334
$ws_BSM.Range("K86").Value = 1615.8334
$ws_BSM.Range("M86").Value = -492.8334

# BSM!row89 (hunk 14: @@ -18931,22 +18940,22 @@)
$ws_BSM.Range("H89").Value = 3449.25
$ws_BSM.Range("I89").Value = 1615.8334
$ws_BSM.Range("K89").Value = 8079.166999999999
$ws_BSM.Range("M89").Value = -2463.166999999999

# BSM!row105 (hunk 15: @@ -19709,25 +19718,25 @@)
$ws_BSM.Range("H105").Value = 8263.823
$ws_BSM.Range("I105").Value = 8320.5
$ws_BSM.Range("J105").Value = 7999.3335
$ws_BSM.Range("K105").Value = 8320.5
$ws_BSM.Range("L105").Value = 7999.3335
$ws_BSM.Range("M105").Value = -6573.5
$ws_BSM.Range("N105").Value = -11493.3335

# CRP!row16 (hunk 16: @@ -22278,25 +22287,22 @@)
$ws_CRP.Range("H16").Value = 23812630
$ws_CRP.Range("I16").Value = 23812630
$ws_CRP.Range("J16").Value = 0
$ws_CRP.Range("K16").Value = 23812630
$ws_CRP.Range("L16").Value = 0
$ws_CRP.Range("M16").Value = -23812343
$ws_CRP.Range("N16").ClearContents()

# CRP!row22 (hunk 17: @@ -22575,22 +22581,22 @@)
$ws_CRP.Range("H22").Value = 1882.6666
$ws_CRP.Range("I22").Value = 1862.909
$ws_CRP.Range("K22").Value = 1862.909
$ws_CRP.Range("M22").Value = -1512.909

# CRP!row105 (hunk 18: @@ -26585,25 +26591,25 @@)
$ws_CRP.Range("H105").Value = 2380.182
$ws_CRP.Range("J105").Value = 4444.3335
$ws_CRP.Range("L105").Value = 4444.3335
$ws_CRP.Range("N105").Value = -7938.3335

# CRP!row113 (hunk 19: @@ -26983,25 +26989,22 @@)
$ws_CRP.Range("H113").Value = 23812630
$ws_CRP.Range("I113").Value = 23812630
$ws_CRP.Range("J113").Value = 0
$ws_CRP.Range("K113").Value = 23812630
$ws_CRP.Range("L113").Value = 0
$ws_CRP.Range("M113").Value = -23810460
$ws_CRP.Range("N113").ClearContents()

# CRP!row122 (hunk 20: @@ -27415,25 +27418,22 @@)
$ws_CRP.Range("H122").Value = 5236
$ws_CRP.Range("I122").Value = 5236
$ws_CRP.Range("J122").Value = 0
$ws_CRP.Range("K122").Value = 15708
$ws_CRP.Range("L122").Value = 0
$ws_CRP.Range("M122").Value = -13258
$ws_CRP.Range("N122").ClearContents()

# CRP!row134 (hunk 21: @@ -27997,22 +27997,22 @@)
$ws_CRP.Range("H134").Value = 1558.8148
$ws_CRP.Range("I134").Value = 1472.7142
$ws_CRP.Range("K134").Value = 4418.142599999999
$ws_CRP.Range("M134").Value = -1883.142599999999

# CUL!row34 (hunk 22: @@ -30129,25 +30129,25 @@)
$ws_CUL.Range("H34").Value = 3089.4666
$ws_CUL.Range("J34").Value = 3514.923
$ws_CUL.Range("L34").Value = 10544.769
$ws_CUL.Range("N34").Value = -10712.769

# CUL!row39 (hunk 23: @@ -30377,25 +30377,25 @@)
$ws_CUL.Range("H39").Value = 6510
$ws_CUL.Range("J39").Value = 14992.5
$ws_CUL.Range("L39").Value = 44977.5
$ws_CUL.Range("N39").Value = -45565.5

# CUL!row55 (hunk 24: @@ -31179,25 +31179,25 @@)
$ws_CUL.Range("H55").Value = 5212.107
$ws_CUL.Range("J55").Value = 5869.1304
$ws_CUL.Range("L55").Value = 17607.3912
$ws_CUL.Range("N55").Value = -17961.3912

# CUL!row138 (hunk 25: @@ -35402,22 +35402,22 @@)
$ws_CUL.Range("H138").Value = 6141.2856
$ws_CUL.Range("I138").Value = 5831.5
$ws_CUL.Range("K138").Value = 17494.5
$ws_CUL.Range("M138").Value = -12354.5

# GSM!row33 (hunk 26: @@ -37223,22 +37223,19 @@)
$ws_GSM.Range("H33").Value = 0
$ws_GSM.Range("I33").Value = 0
$ws_GSM.Range("K33").Value = 0
$ws_GSM.Range("M33").ClearContents()

# GSM!row122 (hunk 27: @@ -41539,25 +41536,25 @@)
$ws_GSM.Range("H122").Value = 2706.4688
$ws_GSM.Range("J122").Value = 3028.4285
$ws_GSM.Range("L122").Value = 9085.2855
$ws_GSM.Range("N122").Value = -13985.2855

# GSM!row126 (hunk 28: @@ -41738,23 +41735,26 @@)
$ws_GSM.Range("H126").Value = 2247.5
$ws_GSM.Range("J126").Value = 2000
$ws_GSM.Range("L126").Value = 6000
$ws_GSM.Range("N126").Value = -10940

# GSM!row132 (hunk 29: @@ -42023,22 +42023,22 @@)
$ws_GSM.Range("H132").Value = 1075708.5
$ws_GSM.Range("I132").Value = 9428.666999999999
$ws_GSM.Range("K132").Value = 28286.001
$ws_GSM.Range("M132").Value = -25756.001

# LTW!row4 (hunk 30: @@ -42705,22 +42705,22 @@)
$ws_LTW.Range("H4").Value = 15047.5
$ws_LTW.Range("I4").Value = 13430.333
$ws_LTW.Range("K4").Value = 13430.333
$ws_LTW.Range("M4").Value = -13317.333

# LTW!row14 (hunk 31: @@ -43195,22 +43195,22 @@)
$ws_LTW.Range("H14").Value = 154950
$ws_LTW.Range("I14").Value = 154950
$ws_LTW.Range("K14").Value = 154950
$ws_LTW.Range("M14").Value = -154778

# LTW!row16 (hunk 32: @@ -43290,22 +43290,22 @@)
$ws_LTW.Range("H16").Value = 808.9524
$ws_LTW.Range("I16").Value = 698.51514
$ws_LTW.Range("K16").Value = 698.51514
$ws_LTW.Range("M16").Value = -528.51514

# LTW!row28 (hunk 33: @@ -43881,22 +43881,22 @@)
$ws_LTW.Range("H28").Value = 15047.5
$ws_LTW.Range("I28").Value = 13430.333
$ws_LTW.Range("K28").Value = 13430.333
$ws_LTW.Range("M28").Value = -13198.333

# LTW!row37 (hunk 34: @@ -44322,22 +44322,22 @@)
$ws_LTW.Range("H37").Value = 15047.5
$ws_LTW.Range("I37").Value = 13430.333
$ws_LTW.Range("K37").Value = 13430.333
$ws_LTW.Range("M37").Value = -13323.333

# LTW!row55 (hunk 35: @@ -45195,22 +45195,22 @@)
$ws_LTW.Range("H55").Value = 1093.2222
$ws_LTW.Range("I55").Value = 947.1429000000001
$ws_LTW.Range("K55").Value = 947.1429000000001
$ws_LTW.Range("M55").Value = -774.1429000000001

# LTW!row122 (hunk 36: @@ -48442,22 +48442,22 @@)
$ws_LTW.Range("H122").Value = 3390.6875
$ws_LTW.Range("I122").Value = 3003.1614
$ws_LTW.Range("K122").Value = 9009.484199999999
$ws_LTW.Range("M122").Value = -6559.484199999999

# LTW!row132 (hunk 37: @@ -48917,25 +48917,25 @@)
$ws_LTW.Range("H132").Value = 3163.182
$ws_LTW.Range("I132").Value = 2791.7144
$ws_LTW.Range("J132").Value = 3813.25
$ws_LTW.Range("K132").Value = 8375.143199999999
$ws_LTW.Range("L132").Value = 11439.75
$ws_LTW.Range("M132").Value = -5845.143199999999
$ws_LTW.Range("N132").Value = -16499.75

# WVR!row122 (hunk 38: @@ -55363,22 +55363,22 @@)
$ws_WVR.Range("H122").Value = 3297.2222
$ws_WVR.Range("I122").Value = 2459.375
$ws_WVR.Range("K122").Value = 7378.125
$ws_WVR.Range("M122").Value = -4928.125

# WVR!row126 (hunk 39: @@ -55562,25 +55562,25 @@)
$ws_WVR.Range("H126").Value = 3030.6924
$ws_WVR.Range("J126").Value = 3486
$ws_WVR.Range("L126").Value = 10458
$ws_WVR.Range("N126").Value = -15398

# WVR!row132 (hunk 40: @@ -55850,25 +55850,25 @@)
$ws_WVR.Range("H132").Value = 2238.0889
$ws_WVR.Range("I132").Value = 1705.1724
$ws_WVR.Range("J132").Value = 3204
$ws_WVR.Range("K132").Value = 5115.5172
$ws_WVR.Range("L132").Value = 9612
$ws_WVR.Range("M132").Value = -2585.5172
$ws_WVR.Range("N132").Value = -14672

# WVR!row135 (hunk 41: @@ -56000,22 +56000,22 @@)
$ws_WVR.Range("H135").Value = 129900
$ws_WVR.Range("J135").Value = 129900
$ws_WVR.Range("L135").Value = 129900
$ws_WVR.Range("N135").Value = -140040

# WVR!row136 (hunk 42: @@ -56049,25 +56049,25 @@)
$ws_WVR.Range("H136").Value = 23287.348
$ws_WVR.Range("J136").Value = 4317.5
$ws_WVR.Range("L136").Value = 12952.5
$ws_WVR.Range("N136").Value = -18052.5

# WVR!row137 (hunk 43: @@ -56101,25 +56101,22 @@)
$ws_WVR.Range("H137").Value = 49500
$ws_WVR.Range("J137").Value = 0
$ws_WVR.Range("L137").Value = 0
$ws_WVR.Range("N137").ClearContents()

Write-Host "Applied scheduled market-data update across 44 leve rows."